# Updates cryptocurrency price/volume figures (and restores two swapped
# rows' Coin/Link values) on the "cryptos" worksheet, per the latest
# scrape from coinranking.com.
#
# Cells in column D hold prices that are stored as *text* (e.g. "333.48",
# "30.777.15") rather than numbers, because some values use "." as a
# thousands separator and would be mangled if Excel auto-detected them as
# numeric. Assigning a plain numeric-looking string via .Value would make
# Excel silently convert the cell to a Number, so for those cells we force
# a Text number format first, assign the value, and then restore the
# cell's original ("Normal") style so no visual/formatting change leaks
# into the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.777.15'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '2.111.67'
$ws.Range('E3').Value = '  +7.32%  '
$ws.Range('E4').Value = '  +0.31%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '333.48'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +3.68%  '
$ws.Range('E6').Value = '  +0.32%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.5312'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  +3.87%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.4397'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +7.48%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.09020'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +7.17%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '46.06'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  +8.42%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '1.178'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +4.35%  '
$ws.Range('E12').Value = '  +3.55%  '
$ws.Range('D13').Value = '2.107.88'
$ws.Range('E13').Value = '  +7.73%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '6.755'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +4.64%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '7.804'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +5.89%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '97.46'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +3.90%  '
$ws.Range('E17').Value = '  +0.23%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '0.00001128'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +2.37%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.06662'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +2.02%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '19.13'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +2.11%  '
$ws.Range('E21').Value = '  +0.17%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '6.356'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +5.23%  '
$ws.Range('D23').Value = '30.833.31'
$ws.Range('E23').Value = '  +2.07%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '12.37'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +7.91%  '
$ws.Range('D25').Value = '2.354.33'
$ws.Range('E25').Value = '  +8.35%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '2.258'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +2.80%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '22.77'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +0.89%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '2.573'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  +9.09%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '162.43'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -0.17%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '132.89'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +2.04%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '1.173'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +2.74%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '0.1078'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +2.32%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '6.226'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +3.24%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '4.018'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +6.09%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '1.543'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  +20.48%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.02605'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +5.43%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '5.534'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +3.36%  '
$ws.Range('E38').Value = '  +3.70%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '9.521'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +7.12%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '12.78'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +8.82%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '0.2275'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +4.93%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.6864'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +4.95%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '1.253'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +2.42%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '14.18'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +4.82%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.6449'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +5.52%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +0.35%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '2.232'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +1.94%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '3.671'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +1.13%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '1.276'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +4.60%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '82.39'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +4.10%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '120.09'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -2.66%  '
